# Generate Report for Handback
#
# This script mirrors a "handback" report refresh: the zh-cn and de-de
# localization status rows move from "Ready for handoff" to
# "Handed back: in sync with en-US", their handback timestamps advance,
# the stale "handback file is not latest" error clears, and the
# "Status"/"Error Detail" columns are widened/narrowed to fit the new
# (longer status / now-empty error) text on all three sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# --- Column width adjustments -------------------------------------------
# (COM ColumnWidth is quantized to the sheet's pixel grid, same as real
# Excel; these inputs land on the pixel step closest to the refreshed
# report's column widths.)

# Overview: widen the "zh-cn" (E) and "de-de" (F) status columns.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn: widen "Status" (C), narrow "Error Detail" (P).
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# de-de: widen "Status" (C), narrow "Error Detail" (P).
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- zh-cn row 2: handback refresh ---------------------------------------
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-09-03 20:51:57"
$zhcn.Range("P2").Value = ""

# --- de-de row 2: handback refresh ---------------------------------------
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-09-03 20:52:10"
$dede.Range("P2").Value = ""
